# "Changed All Tests to Read from Excel Sheet"
# Rename the existing sheet, add two more sheets with test data, and update
# selections to match the new authoring session.

$wb = $excel.ActiveWorkbook

# --- Sheet1: Sheet1 -> LoginTest -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "LoginTest"

# --- Sheet2: NewMedicineService ---------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "NewMedicineService"

# Header row (bold, Text number format)
$ws2.Range("A1").Value = "TestCaseName"
$ws2.Range("B1").Value = "ProductSearch"
$ws2.Range("C1").Value = "inProgressClinicalServicesExpected"
$ws2.Range("D1").Value = "InterventionDate"
$ws2.Range("E1").Value = "FollowUpDate"
$ws2.Range("F1").Value = "deliveredClinicalServicesExpected"
$ws2.Rows.Item(1).Font.Bold = $true
$ws2.Range("A1:F1").NumberFormat = "@"

# Data row -- text cells get the Text format applied before the value is
# entered (so numeric-looking text like "30" stays text); the numeric cell
# (C2) gets the value first and the format applied after (stays a real
# number formatted as text).
$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "Warfarin"

$ws2.Range("D2:E2").NumberFormat = "@"
$ws2.Range("D2").Value = "13/01/2021"
$ws2.Range("E2").Value = "13/01/2021"

$ws2.Range("F2").NumberFormat = "@"
$ws2.Range("F2").Value = "30"

$ws2.Range("A2").NumberFormat = "@"
$ws2.Range("A2").Value = "NewMedicine1"

$ws2.Range("C2").Value = 36
$ws2.Range("C2").NumberFormat = "@"

# Column widths (best effort -- engine quantizes to 1/6 character units)
$ws2.Columns.Item(1).ColumnWidth = 15.665
$ws2.Columns.Item(2).ColumnWidth = 15.165
$ws2.Columns.Item(3).ColumnWidth = 30.33
$ws2.Columns.Item(4).ColumnWidth = 16.33
$ws2.Columns.Item(5).ColumnWidth = 16.83
$ws2.Columns.Item(6).ColumnWidth = 30.0

$ws2.Range("D18").Select() | Out-Null

# --- Sheet3: PatientRecords -------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "PatientRecords"

$ws3.Range("A1").Value = "TestCaseName"
$ws3.Range("B1").Value = "Surname"
$ws3.Range("C1").Value = "FirstName"
$ws3.Rows.Item(1).Font.Bold = $true

$ws3.Range("A2").Value = "PatientSearch1"
$ws3.Range("B2").Value = "Smith"
$ws3.Range("C2").Value = "Robert"

$ws3.Columns.Item(1).ColumnWidth = 14.5
$ws3.Columns.Item(2).ColumnWidth = 12.83
$ws3.Columns.Item(3).ColumnWidth = 13.665

$ws3.Range("C1").Select() | Out-Null

# --- Back on Sheet1: move the selection, keep it the active tab ------------------
$ws1.Range("D5").Select() | Out-Null
$ws1.Activate() | Out-Null
